$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = $fmt
}

$wsCPB = $wb.Worksheets.Item("CPB")
$wsTRP = $wb.Worksheets.Item("TRP")

Set-TextValue $wsCPB "D6" "13194"
Set-TextValue $wsCPB "E6" "13624"
Set-TextValue $wsCPB "D7" "88102"
Set-TextValue $wsCPB "E7" "90973"
Set-TextValue $wsCPB "D9" "6078"
Set-TextValue $wsCPB "E9" "6276"
Set-TextValue $wsCPB "D10" "7881"
Set-TextValue $wsCPB "E10" "8138"
Set-TextValue $wsCPB "D11" "81616"
Set-TextValue $wsCPB "E11" "84276"
Set-TextValue $wsCPB "D12" "20507"
Set-TextValue $wsCPB "E12" "21175"
Set-TextValue $wsCPB "D14" "11007"
Set-TextValue $wsCPB "E14" "11366"
Set-TextValue $wsCPB "D15" "11234"
Set-TextValue $wsCPB "E15" "11600"
Set-TextValue $wsCPB "D16" "5509"
Set-TextValue $wsCPB "E16" "5689"
Set-TextValue $wsCPB "D17" "6347"
Set-TextValue $wsCPB "E17" "6554"
Set-TextValue $wsCPB "D19" "5975"
Set-TextValue $wsCPB "E19" "6170"
Set-TextValue $wsCPB "D20" "5088"
Set-TextValue $wsCPB "E20" "5254"
Set-TextValue $wsCPB "D21" "27343"
Set-TextValue $wsCPB "E21" "28234"
Set-TextValue $wsCPB "D22" "6059"
Set-TextValue $wsCPB "E22" "6256"
Set-TextValue $wsCPB "D24" "4250"
Set-TextValue $wsCPB "E24" "4389"
Set-TextValue $wsCPB "D25" "10804"
Set-TextValue $wsCPB "E25" "11156"
Set-TextValue $wsCPB "D27" "5116"
Set-TextValue $wsCPB "E27" "5283"
Set-TextValue $wsCPB "D28" "4260"
Set-TextValue $wsCPB "E28" "4399"
Set-TextValue $wsCPB "D35" "3004"
Set-TextValue $wsCPB "E35" "3102"
Set-TextValue $wsCPB "D37" "1934"
Set-TextValue $wsCPB "E37" "1997"
Set-TextValue $wsCPB "D40" "4564"
Set-TextValue $wsCPB "E40" "4713"
Set-TextValue $wsTRP "D6" "13424"
Set-TextValue $wsTRP "E6" "13861,5"
Set-TextValue $wsTRP "D7" "120808"
Set-TextValue $wsTRP "E7" "124745,5"
Set-TextValue $wsTRP "D8" "8797"
Set-TextValue $wsTRP "E8" "9083,7"
Set-TextValue $wsTRP "D9" "6165"
Set-TextValue $wsTRP "E9" "6365,9"
Set-TextValue $wsTRP "D10" "7938"
Set-TextValue $wsTRP "E10" "8196,7"
Set-TextValue $wsTRP "D11" "82617"
Set-TextValue $wsTRP "E11" "85309,7"
Set-TextValue $wsTRP "D12" "20672"
Set-TextValue $wsTRP "E12" "21345,8"
Set-TextValue $wsTRP "D13" "5729"
Set-TextValue $wsTRP "E13" "5915,7"
Set-TextValue $wsTRP "D14" "11088"
Set-TextValue $wsTRP "E14" "11449,4"
Set-TextValue $wsTRP "D15" "11322"
Set-TextValue $wsTRP "E15" "11691"
Set-TextValue $wsTRP "D16" "5587"
Set-TextValue $wsTRP "E16" "5769,1"
Set-TextValue $wsTRP "D17" "6410"
Set-TextValue $wsTRP "E17" "6618,9"
Set-TextValue $wsTRP "D18" "6173"
Set-TextValue $wsTRP "E18" "6374,2"
Set-TextValue $wsTRP "D19" "6025"
Set-TextValue $wsTRP "E19" "6221,4"
Set-TextValue $wsTRP "E20" "5273,5"
Set-TextValue $wsTRP "D21" "27299"
Set-TextValue $wsTRP "E21" "28188,8"
Set-TextValue $wsTRP "D22" "5409"
Set-TextValue $wsTRP "E22" "5585,3"
Set-TextValue $wsTRP "E23" "5927,1"
Set-TextValue $wsTRP "E24" "4380,3"
Set-TextValue $wsTRP "D25" "10800"
Set-TextValue $wsTRP "E25" "11152"
Set-TextValue $wsTRP "E26" "4605,4"
Set-TextValue $wsTRP "E27" "5224,9"
Set-TextValue $wsTRP "D28" "4256"
Set-TextValue $wsTRP "E28" "4394,7"
Set-TextValue $wsTRP "E29" "2610,4"
Set-TextValue $wsTRP "E30" "2670,3"
Set-TextValue $wsTRP "D31" "2605"
Set-TextValue $wsTRP "E31" "2689,9"
Set-TextValue $wsTRP "D32" "2740"
Set-TextValue $wsTRP "E32" "2829,3"
Set-TextValue $wsTRP "D33" "3612"
Set-TextValue $wsTRP "E33" "3729,7"
Set-TextValue $wsTRP "D34" "4032"
Set-TextValue $wsTRP "E34" "4163,4"
Set-TextValue $wsTRP "E36" "3415,8"
Set-TextValue $wsTRP "E37" "2003,2"
Set-TextValue $wsTRP "E38" "2105,5"
Set-TextValue $wsTRP "E39" "1722,4"
Set-TextValue $wsTRP "E40" "4701,4"
Set-TextValue $wsTRP "E41" "3077,1"
Set-TextValue $wsTRP "D43" "8070"
Set-TextValue $wsTRP "E43" "8333"
